$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.106.60"
$ws.Range("E2").Value = "  -1.34%  "

# Row 3
$ws.Range("D3").Value = "3.383.58"
$ws.Range("E3").Value = "  -0.34%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.77"
$ws.Range("E5").Value = "  -0.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.29"
$ws.Range("E6").Value = "  -0.81%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").Value = "3.382.14"
$ws.Range("E8").Value = "  -0.39%  "

# Row 9
$ws.Range("E9").Value = "  -1.06%  "

# Row 10
$ws.Range("E10").Value = "  +2.35%  "

# Row 11
$ws.Range("E11").Value = "  -3.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.382"
$ws.Range("E12").Value = "  -2.62%  "

# Row 13
$ws.Range("D13").Value = "3.962.21"
$ws.Range("E13").Value = "  -0.33%  "

# Row 14
$ws.Range("E14").Value = "  +0.94%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  -3.42%  "

# Row 16
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.65"
$ws.Range("E16").Value = "  +0.83%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.390.63"
$ws.Range("E17").Value = "  -0.19%  "

# Row 18
$ws.Range("D18").Value = "61.266.06"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.85"
$ws.Range("E19").Value = "  -2.22%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.36"
$ws.Range("E20").Value = "  -1.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.74"
$ws.Range("E21").Value = "  -1.33%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.69"
$ws.Range("E22").Value = "  -4.21%  "

# Row 23
$ws.Range("B23").Value = "WrappedeETH"
$ws.Range("C23").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D23").Value = "3.525.20"
$ws.Range("E23").Value = "  -0.55%  "

# Row 24
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.552"
$ws.Range("E24").Value = "  -2.50%  "

# Row 25
$ws.Range("E25").Value = "  +0.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  -2.60%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.12"
$ws.Range("E27").Value = "  -0.69%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.183"
$ws.Range("E28").Value = "  +13.12%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.64"
$ws.Range("E29").Value = "  -0.96%  "

# Row 30
$ws.Range("E30").Value = "  -0.27%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.38"
$ws.Range("E31").Value = "  -3.47%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  -1.64%  "

# Row 33
$ws.Range("E33").Value = "  -1.58%  "

# Row 34
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.44"
$ws.Range("E35").Value = "  -0.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.16"
$ws.Range("E36").Value = "  -4.59%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  -3.36%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.83"
$ws.Range("E38").Value = "  -1.20%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.41"
$ws.Range("E39").Value = "  -0.32%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0758"
$ws.Range("E40").Value = "  -3.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.52"
$ws.Range("E41").Value = "  +2.22%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.04%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.773"
$ws.Range("E43").Value = "  -1.79%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.70"
$ws.Range("E44").Value = "  +0.83%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.70"
$ws.Range("E45").Value = "  -3.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.35"
$ws.Range("E46").Value = "  -1.95%  "

# Row 47
$ws.Range("E47").Value = "  -5.30%  "

# Row 48
$ws.Range("D48").Value = "2.514.53"
$ws.Range("E48").Value = "  +7.55%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.78"
$ws.Range("E49").Value = "  -1.81%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.92"
$ws.Range("E50").Value = "  -0.78%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.46"
$ws.Range("E51").Value = "  +4.68%  "
